$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and E stay text cells (matches original inlineStr storage),
# so numeric-looking values like '7.00' or '2.00' keep their trailing zeros.
$ws.Range('D2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '57.974.39'
$ws.Range('E2').Value = '  +1.51%  '
$ws.Range('D3').Value = '2.943.82'
$ws.Range('E3').Value = '  +3.22%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '554.38'
$ws.Range('E5').Value = '  +2.46%  '
$ws.Range('D6').Value = '133.26'
$ws.Range('E6').Value = '  +11.21%  '
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('D8').Value = '0.512'
$ws.Range('E8').Value = '  +6.48%  '
$ws.Range('D9').Value = '2.941.66'
$ws.Range('E9').Value = '  +3.37%  '
$ws.Range('E10').Value = '  +4.48%  '
$ws.Range('D11').Value = '4.79'
$ws.Range('E11').Value = '  +1.56%  '
$ws.Range('D12').Value = '0.447'
$ws.Range('E12').Value = '  +6.33%  '
$ws.Range('D13').Value = '0.0000221'
$ws.Range('E13').Value = '  +6.44%  '
$ws.Range('D14').Value = '32.92'
$ws.Range('E14').Value = '  +7.10%  '
$ws.Range('E15').Value = '  +3.83%  '
$ws.Range('D16').Value = '3.428.17'
$ws.Range('E16').Value = '  +3.36%  '
$ws.Range('D17').Value = '6.91'
$ws.Range('E17').Value = '  +12.30%  '
$ws.Range('D18').Value = '2.933.93'
$ws.Range('E18').Value = '  +3.42%  '
$ws.Range('D19').Value = '57.935.03'
$ws.Range('E19').Value = '  +1.39%  '
$ws.Range('D20').Value = '417.33'
$ws.Range('E20').Value = '  +3.04%  '
$ws.Range('D21').Value = '13.32'
$ws.Range('E21').Value = '  +6.65%  '
$ws.Range('D22').Value = '0.695'
$ws.Range('E22').Value = '  +9.52%  '
$ws.Range('D23').Value = '13.42'
$ws.Range('E23').Value = '  +9.52%  '
$ws.Range('D24').Value = '7.00'
$ws.Range('E24').Value = '  +5.84%  '
$ws.Range('D25').Value = '78.98'
$ws.Range('E25').Value = '  +4.90%  '
$ws.Range('E26').Value = '  -0.09%  '
$ws.Range('D27').Value = '0.998'
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('D28').Value = '2.49'
$ws.Range('E28').Value = '  +4.23%  '
$ws.Range('D29').Value = '2.03'
$ws.Range('E29').Value = '  +8.31%  '
$ws.Range('D30').Value = '7.47'
$ws.Range('E30').Value = '  +6.95%  '
$ws.Range('D31').Value = '25.50'
$ws.Range('E31').Value = '  +5.08%  '
$ws.Range('D32').Value = '5.95'
$ws.Range('E32').Value = '  +2.90%  '
$ws.Range('D33').Value = '0.0977'
$ws.Range('E33').Value = '  +7.43%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = '5.69'
$ws.Range('E34').Value = '  +8.18%  '
$ws.Range('B35').Value = 'Mantle'
$ws.Range('C35').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D35').Value = '0.950'
$ws.Range('E35').Value = '  +8.51%  '
$ws.Range('D36').Value = '2.08'
$ws.Range('E36').Value = '  +6.17%  '
$ws.Range('D37').Value = '0.0₃0703'
$ws.Range('E37').Value = '  +16.13%  '
$ws.Range('D38').Value = '48.24'
$ws.Range('E38').Value = '  +0.24%  '
$ws.Range('D39').Value = '8.75'
$ws.Range('E39').Value = '  +7.52%  '
$ws.Range('D40').Value = '2.66'
$ws.Range('E40').Value = '  +15.27%  '
$ws.Range('D41').Value = '380.62'
$ws.Range('E41').Value = '  +10.80%  '
$ws.Range('E42').Value = '  +5.73%  '
$ws.Range('D43').Value = '0.0348'
$ws.Range('E43').Value = '  +4.24%  '
$ws.Range('D44').Value = '2.701.06'
$ws.Range('E44').Value = '  +5.34%  '
$ws.Range('E45').Value = '  +0.07%  '
$ws.Range('D46').Value = '124.36'
$ws.Range('E46').Value = '  +6.58%  '
$ws.Range('D47').Value = '0.237'
$ws.Range('E47').Value = '  +6.48%  '
$ws.Range('B48').Value = 'Fetch.AI'
$ws.Range('C48').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D48').Value = '1.97'
$ws.Range('E48').Value = '  +4.61%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').Value = '0.107'
$ws.Range('E49').Value = '  +3.56%  '
$ws.Range('D50').Value = '22.93'
$ws.Range('E50').Value = '  +3.67%  '
$ws.Range('D51').Value = '2.00'
$ws.Range('E51').Value = '  +5.98%  '
